$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "AI, 일자리 공포를 넘어 ‘두려움 없는 미래’로!"
$ws.Range("E18").Value = "https://freesearch.pe.kr/archives/5526"

$ws.Range("D23").Value = "Evaluate expression to the main debugger toolbar (new ui)"
$ws.Range("E23").Value = "https://theonly1.tistory.com/3326"

$ws.Range("D27").Value = "핑퐁"
